$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
# Row 51
$ws.Range("H51").Value = 2600
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2600
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2600
$ws.Range("M51").Value = 0
$ws.Range("N51").Value = -3568
# Row 70
$ws.Range("H70").Value = 1586.25
$ws.Range("I70").Value = 1390
$ws.Range("J70").Value = 1614.2858
$ws.Range("K70").Value = 4170
$ws.Range("L70").Value = 4842.857400000001
$ws.Range("M70").Value = -3900
$ws.Range("N70").Value = -5382.857400000001
# Row 73
$ws.Range("H73").Value = 1586.25
$ws.Range("I73").Value = 1390
$ws.Range("J73").Value = 1614.2858
$ws.Range("K73").Value = 4170
$ws.Range("L73").Value = 4842.857400000001
$ws.Range("M73").Value = -3234
$ws.Range("N73").Value = -6714.857400000001
# Row 107
$ws.Range("H107").Value = 12334.889
$ws.Range("I107").Value = 12334.889
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 12334.889
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -10414.889
$ws.Range("N107").ClearContents()
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 137
$ws.Range("H137").Value = 1015.6667
$ws.Range("I137").Value = 859.5
$ws.Range("J137").Value = 1562.25
$ws.Range("K137").Value = 2578.5
$ws.Range("L137").Value = 4686.75
$ws.Range("M137").Value = -28.5
$ws.Range("N137").Value = -9786.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 321356.16
$ws.Range("I32").Value = 4962.0156
$ws.Range("J32").Value = 1071327.5
$ws.Range("K32").Value = 4962.0156
$ws.Range("L32").Value = 1071327.5
$ws.Range("M32").Value = -4675.0156
$ws.Range("N32").Value = -1071901.5
# Row 45
$ws.Range("H45").Value = 2951.4666
$ws.Range("I45").Value = 2782.0557
$ws.Range("J45").Value = 3205.5833
$ws.Range("K45").Value = 2782.0557
$ws.Range("L45").Value = 3205.5833
$ws.Range("M45").Value = -2405.0557
$ws.Range("N45").Value = -3959.5833
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()
# Row 132
$ws.Range("H132").Value = 2182.9285
$ws.Range("I132").Value = 1876.6818
$ws.Range("K132").Value = 5630.0454
$ws.Range("M132").Value = -3100.0454

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1365.2391
$ws.Range("I134").Value = 838.9091
$ws.Range("J134").Value = 2701.3076
$ws.Range("K134").Value = 2516.7273
$ws.Range("L134").Value = 8103.9228
$ws.Range("M134").Value = 18.27269999999999
$ws.Range("N134").Value = -13173.9228

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 2322.2258
$ws.Range("I99").Value = 1860.9565
$ws.Range("J99").Value = 3648.375
$ws.Range("K99").Value = 1860.9565
$ws.Range("L99").Value = 3648.375
$ws.Range("M99").Value = -362.9565
$ws.Range("N99").Value = -6644.375
# Row 105
$ws.Range("H105").Value = 1114.0834
$ws.Range("I105").Value = 763.3333
$ws.Range("J105").Value = 2166.3333
$ws.Range("K105").Value = 763.3333
$ws.Range("L105").Value = 2166.3333
$ws.Range("M105").Value = 983.6667
$ws.Range("N105").Value = -5660.3333
# Row 122
$ws.Range("H122").Value = 2508
$ws.Range("I122").Value = 1666.6666
$ws.Range("J122").Value = 3139
$ws.Range("K122").Value = 4999.9998
$ws.Range("L122").Value = 9417
$ws.Range("M122").Value = -2549.9998
$ws.Range("N122").Value = -14317
# Row 126
$ws.Range("H126").Value = 2322.2258
$ws.Range("I126").Value = 1860.9565
$ws.Range("J126").Value = 3648.375
$ws.Range("K126").Value = 5582.8695
$ws.Range("L126").Value = 10945.125
$ws.Range("M126").Value = -3112.8695
$ws.Range("N126").Value = -15885.125
# Row 132
$ws.Range("H132").Value = 2738.238
$ws.Range("I132").Value = 1815.7693
$ws.Range("J132").Value = 4237.25
$ws.Range("K132").Value = 5447.3079
$ws.Range("L132").Value = 12711.75
$ws.Range("M132").Value = -2917.3079
$ws.Range("N132").Value = -17771.75
# Row 134
$ws.Range("H134").Value = 835.1818
$ws.Range("I134").Value = 661
$ws.Range("J134").Value = 1810.6
$ws.Range("K134").Value = 1983
$ws.Range("L134").Value = 5431.799999999999
$ws.Range("M134").Value = 552
$ws.Range("N134").Value = -10501.8

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 135
$ws.Range("I7").Value = 90
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 270
$ws.Range("L7").Value = 450
$ws.Range("M7").Value = -158
$ws.Range("N7").Value = -674

$ws = $wb.Worksheets.Item("GSM")
# Row 99
$ws.Range("H99").Value = 9166.666999999999
$ws.Range("I99").Value = 9166.666999999999
$ws.Range("K99").Value = 9166.666999999999
$ws.Range("M99").Value = -6920.666999999999
# Row 122
$ws.Range("H122").Value = 1998
$ws.Range("I122").Value = 1818.2858
$ws.Range("J122").Value = 2627
$ws.Range("K122").Value = 5454.857400000001
$ws.Range("L122").Value = 7881
$ws.Range("M122").Value = -3004.857400000001
$ws.Range("N122").Value = -12781
# Row 132
$ws.Range("H132").Value = 1910.9412
$ws.Range("I132").Value = 1500.6111
$ws.Range("K132").Value = 4501.8333
$ws.Range("M132").Value = -1971.8333
# Row 136
$ws.Range("H136").Value = 21262.584
$ws.Range("J136").Value = 21262.584
$ws.Range("L136").Value = 63787.75199999999
$ws.Range("N136").Value = -68887.75199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5791.5713
$ws.Range("I62").Value = 6497.4287
$ws.Range("J62").Value = 5085.7144
$ws.Range("K62").Value = 6497.4287
$ws.Range("L62").Value = 5085.7144
$ws.Range("M62").Value = -5873.4287
$ws.Range("N62").Value = -6333.7144
# Row 65
$ws.Range("H65").Value = 5791.5713
$ws.Range("I65").Value = 6497.4287
$ws.Range("J65").Value = 5085.7144
$ws.Range("K65").Value = 32487.1435
$ws.Range("L65").Value = 25428.572
$ws.Range("M65").Value = -29367.1435
$ws.Range("N65").Value = -31668.572
# Row 132
$ws.Range("H132").Value = 35715904
$ws.Range("I132").Value = 55556930
$ws.Range("J132").Value = 2060.0667
$ws.Range("K132").Value = 166670790
$ws.Range("L132").Value = 6180.2001
$ws.Range("M132").Value = -166668260
$ws.Range("N132").Value = -11240.2001
# Row 136
$ws.Range("H136").Value = 658.3684
$ws.Range("I136").Value = 635.32654
$ws.Range("J136").Value = 799.5
$ws.Range("K136").Value = 1905.97962
$ws.Range("L136").Value = 2398.5
$ws.Range("M136").Value = 644.0203799999999
$ws.Range("N136").Value = -7498.5
